$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Order Number column keeps its leading zeros by forcing text format
$orderRange = $ws.Range("D2:D18")
$orderRange.NumberFormat = "@"
$orderRange.Value = "000002745"

# Update the Time column (G) for rows 2-18
$ws.Range("G2:G18").Value = "23:05:19 2024-05-21"
